$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had columns A:G (dea_code, update_month, month_actual,
# month_forecast, ytd_forecast, year_forecast, loa_forecast). A new column
# "year_remaining_balance" is being inserted before the existing
# "loa_forecast" column, so loa_forecast moves from G to H and the new
# header takes its place in G.

# Remember the current width of column G (loa_forecast) before we touch it,
# so the moved column can keep its original best-fit size.
$locForecastWidth = $ws.Columns.Item(7).ColumnWidth
$locForecastHeader = $ws.Range("G1").Text

# Move loa_forecast header to the new column H
$ws.Range("H1").Value = $locForecastHeader

# Put the new header in column G
$ws.Range("G1").Value = "year_remaining_balance"

# Resize the columns to (best) fit their new contents
$ws.Columns.Item(7).ColumnWidth = 20.67
$ws.Columns.Item(8).ColumnWidth = $locForecastWidth

# Update the active selection like the source workbook shows
$ws.Range("G5").Select()

# Record a page setup (portrait) as reflected in the saved workbook
$ws.PageSetup.Orientation = 1
